$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "tap water production, underground water with disinfection" activity
# was replaced by "market for water, deionised" in the sensitivity table.
$ws.Range("A14").Value = "market for water, deionised"

# The sensitivity values in column B were recalculated against the new
# upstream process; update them in place while keeping them stored as text
# (as they were originally) rather than letting Excel auto-convert them to
# numbers.
$valsRange = $ws.Range("B2:B18")
$valsRange.NumberFormat = "@"

$ws.Range("B2").Value = "1.0915011789959204"
$ws.Range("B3").Value = "1.4044144667855358"
$ws.Range("B4").Value = "2.896388497525361"
$ws.Range("B5").Value = "7.072326777651805"
$ws.Range("B6").Value = "28.157543489605036"
$ws.Range("B7").Value = "3.608264642782406"
$ws.Range("B8").Value = "7.3391499303980705"
$ws.Range("B9").Value = "0.503052879013989"
$ws.Range("B10").Value = "7.548780280013951"
$ws.Range("B11").Value = "2.2490876581265202"
$ws.Range("B12").Value = "0.520341303464715"
$ws.Range("B13").Value = "0.0044123921630157305"
$ws.Range("B14").Value = "0.02796440517023021"
$ws.Range("B15").Value = "3.5374012957409255"
$ws.Range("B16").Value = "0.005401117245006083"
$ws.Range("B17").Value = "0.15758697085276432"
$ws.Range("B18").Value = "-0.25951011368274396"

# Restore the original (default) cell style now that the text values are
# locked in, so formatting matches the source workbook.
$valsRange.Style = "Normal"
